{"js": "// Update the worksheet date title and every \"A\u00d7B=C\" answer cell in the\n// practice table to the next day's values (2024-06-13 -> 2024-06-14).\n// Each data row/col pair is addressed directly via Table.getCell(row, col)\n// (physical table-row index, 0-based) so the edit is unambiguous even\n// though some \"before\"/\"after\" strings repeat across cells.\n\nconst body = context.document.body;\n\n// --- Title paragraph: \"2024-06-13 Thursday\" -> \"2024-06-14 Friday\" ------\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\nparagraphs.items[0].insertText(\"2024-06-14 Friday\", \"Replace\");\n\n// --- Table answer cells ---------------------------------------------------\nconst table = body.tables.getFirst();\n\n// [tableRowIndex, columnIndex, newValue]\nconst updates = [\n  [0, 0, \"510\u00d77=3570\"],\n  [0, 1, \"644\u00d73=1932\"],\n  [0, 2, \"773\u00d73=2319\"],\n  [0, 3, \"524\u00d76=3144\"],\n  [0, 4, \"224\u00d78=1792\"],\n\n  [4, 0, \"365\u00d76=2190\"],\n  [4, 1, \"527\u00d77=3689\"],\n  [4, 2, \"793\u00d77=5551\"],\n  [4, 3, \"661\u00d79=5949\"],\n  [4, 4, \"407\u00d76=2442\"],\n\n  [9, 0, \"388\u00d73=1164\"],\n  [9, 1, \"268\u00d72=536\"],\n  [9, 2, \"259\u00d76=1554\"],\n  [9, 3, \"453\u00d75=2265\"],\n  [9, 4, \"938\u00d78=7504\"],\n\n  [14, 0, \"490\u00d77=3430\"],\n  [14, 1, \"350\u00d79=3150\"],\n  [14, 2, \"661\u00d79=5949\"],\n  [14, 3, \"907\u00d72=1814\"],\n  [14, 4, \"420\u00d74=1680\"],\n\n  [19, 0, \"263\u00d76=1578\"],\n  [19, 1, \"949\u00d73=2847\"],\n  [19, 2, \"463\u00d78=3704\"],\n  [19, 3, \"649\u00d78=5192\"],\n  [19, 4, \"781\u00d78=6248\"],\n];\n\nfor (const [row, col, value] of updates) {\n  table.getCell(row, col).value = value;\n}\n\nawait context.sync();\n", "ps1": "# Update the worksheet date title and every \"A\u00d7B=C\" answer cell in the\n# practice table to the next day's values (2024-06-13 -> 2024-06-14).\n# Each data row/col pair is addressed directly via Table.Cell(row, col)\n# (1-based, matching Word COM conventions) so the edit is unambiguous even\n# though some \"before\"/\"after\" strings repeat across cells.\n\n$d = $word.ActiveDocument\n\n# --- Title paragraph: \"2024-06-13 Thursday\" -> \"2024-06-14 Friday\" ------\n$d.Paragraphs.Item(1).Range.Text = \"2024-06-14 Friday\"\n\n# --- Table answer cells ---------------------------------------------------\n$t = $d.Tables.Item(1)\n\n$t.Cell(1, 1).Range.Text = \"510\u00d77=3570\"\n$t.Cell(1, 2).Range.Text = \"644\u00d73=1932\"\n$t.Cell(1, 3).Range.Text = \"773\u00d73=2319\"\n$t.Cell(1, 4).Range.Text = \"524\u00d76=3144\"\n$t.Cell(1, 5).Range.Text = \"224\u00d78=1792\"\n\n$t.Cell(5, 1).Range.Text = \"365\u00d76=2190\"\n$t.Cell(5, 2).Range.Text = \"527\u00d77=3689\"\n$t.Cell(5, 3).Range.Text = \"793\u00d77=5551\"\n$t.Cell(5, 4).Range.Text = \"661\u00d79=5949\"\n$t.Cell(5, 5).Range.Text = \"407\u00d76=2442\"\n\n$t.Cell(10, 1).Range.Text = \"388\u00d73=1164\"\n$t.Cell(10, 2).Range.Text = \"268\u00d72=536\"\n$t.Cell(10, 3).Range.Text = \"259\u00d76=1554\"\n$t.Cell(10, 4).Range.Text = \"453\u00d75=2265\"\n$t.Cell(10, 5).Range.Text = \"938\u00d78=7504\"\n\n$t.Cell(15, 1).Range.Text = \"490\u00d77=3430\"\n$t.Cell(15, 2).Range.Text = \"350\u00d79=3150\"\n$t.Cell(15, 3).Range.Text = \"661\u00d79=5949\"\n$t.Cell(15, 4).Range.Text = \"907\u00d72=1814\"\n$t.Cell(15, 5).Range.Text = \"420\u00d74=1680\"\n\n$t.Cell(20, 1).Range.Text = \"263\u00d76=1578\"\n$t.Cell(20, 2).Range.Text = \"949\u00d73=2847\"\n$t.Cell(20, 3).Range.Text = \"463\u00d78=3704\"\n$t.Cell(20, 4).Range.Text = \"649\u00d78=5192\"\n$t.Cell(20, 5).Range.Text = \"781\u00d78=6248\"\n"}
